$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rebuild the LR-pair table (A2:T21) with the updated TPM-derived values.
# New data is a full cross-product of 4 sending clusters x 5 target clusters
# (the previous version omitted "Inflammatory-Mac" as a target cluster).
$data = @(
    @('ECs', 'Cadm1', 'Cadm3', 'ECs', 3, 1, 0.5396133333333334, 1.61884, 0.02587018426425635, 0.02587018426425635, 3, 1, 8.725531, 26.176593, 0.6695022935100516, 0.6695022935100517, 4.708412868013334, 42.37571581212, 0.01732014769844728, 0.01732014769844728),
    @('ECs', 'Cadm1', 'Cadm3', 'FAPs', 3, 1, 0.5396133333333334, 1.61884, 0.02587018426425635, 0.02587018426425635, 3, 1, 2.676334666666667, 8.029004, 0.2053527971574215, 0.2053527971574215, 1.444185870595556, 12.99767283536, 0.005312514701642953, 0.005312514701642953),
    @('ECs', 'Cadm1', 'Cadm3', 'Inflammatory-Mac', 3, 1, 0.5396133333333334, 1.61884, 0.02587018426425635, 0.02587018426425635, 3, 1, 0.735084, 2.205252, 0.05640234662194688, 0.05640234662194689, 0.39666112752, 3.56995014768, 0.001459139100046223, 0.001459139100046223),
    @('ECs', 'Cadm1', 'Cadm3', 'MuSCs', 3, 1, 0.5396133333333334, 1.61884, 0.02587018426425635, 0.02587018426425635, 3, 1, 0.8813406666666667, 2.644022, 0.06762449158647323, 0.06762449158647324, 0.4755831749422222, 4.280248574480001, 0.001749458058118716, 0.001749458058118716),
    @('ECs', 'Cadm1', 'Cadm3', 'Resolving-Mac', 3, 1, 0.5396133333333334, 1.61884, 0.02587018426425635, 0.02587018426425635, 1, 0.3333333333333333, 0.01457166666666667, 0.043715, 0.001118071124106636, 0.001118071124106636, 0.007863065622222222, 0.0707675906, 0.0000289247060011829094353, 0.0000289247060011829195997),
    @('Inflammatory-Mac', 'Cadm1', 'Cadm3', 'ECs', 3, 1, 6.260434, 18.781302, 0.3001382122153186, 0.3001382122153186, 3, 1, 8.725531, 26.176593, 0.6695022935100516, 0.6695022935100517, 54.625610940454, 491.630498464086, 0.2009432214481624, 0.2009432214481624),
    @('Inflammatory-Mac', 'Cadm1', 'Cadm3', 'FAPs', 3, 1, 6.260434, 18.781302, 0.3001382122153186, 0.3001382122153186, 3, 1, 2.676334666666667, 8.029004, 0.2053527971574215, 0.2053527971574215, 16.75501654257867, 150.795148883208, 0.06163422141224345, 0.06163422141224346),
    @('Inflammatory-Mac', 'Cadm1', 'Cadm3', 'Inflammatory-Mac', 3, 1, 6.260434, 18.781302, 0.3001382122153186, 0.3001382122153186, 3, 1, 0.735084, 2.205252, 0.05640234662194688, 0.05640234662194689, 4.601944866456, 41.41750379810399, 0.01692849947985985, 0.01692849947985985),
    @('Inflammatory-Mac', 'Cadm1', 'Cadm3', 'MuSCs', 3, 1, 6.260434, 18.781302, 0.3001382122153186, 0.3001382122153186, 3, 1, 0.8813406666666667, 2.644022, 0.06762449158647323, 0.06762449158647324, 5.517575075182667, 49.658175676644, 0.02029669400673393, 0.02029669400673393),
    @('Inflammatory-Mac', 'Cadm1', 'Cadm3', 'Resolving-Mac', 3, 1, 6.260434, 18.781302, 0.3001382122153186, 0.3001382122153186, 1, 0.3333333333333333, 0.01457166666666667, 0.043715, 0.001118071124106636, 0.001118071124106636, 0.09122495743666666, 0.82102461693, 0.0003355758683189374, 0.0003355758683189375),
    @('MuSCs', 'Cadm1', 'Cadm3', 'ECs', 3, 1, 1.182122, 3.546366, 0.0566733845769154, 0.0566733845769154, 3, 1, 8.725531, 26.176593, 0.6695022935100516, 0.6695022935100517, 10.314642156782, 92.831779411038, 0.03794296095522205, 0.03794296095522206),
    @('MuSCs', 'Cadm1', 'Cadm3', 'FAPs', 3, 1, 1.182122, 3.546366, 0.0566733845769154, 0.0566733845769154, 3, 1, 2.676334666666667, 8.029004, 0.2053527971574215, 0.2053527971574215, 3.163754088829334, 28.473786799464, 0.01163803804724785, 0.01163803804724785),
    @('MuSCs', 'Cadm1', 'Cadm3', 'Inflammatory-Mac', 3, 1, 1.182122, 3.546366, 0.0566733845769154, 0.0566733845769154, 3, 1, 0.735084, 2.205252, 0.05640234662194688, 0.05640234662194689, 0.8689589682479999, 7.820630714231999, 0.003196511881146081, 0.003196511881146081),
    @('MuSCs', 'Cadm1', 'Cadm3', 'MuSCs', 3, 1, 1.182122, 3.546366, 0.0566733845769154, 0.0566733845769154, 3, 1, 0.8813406666666667, 2.644022, 0.06762449158647323, 0.06762449158647324, 1.041852191561333, 9.376669724052, 0.003832508818498577, 0.003832508818498578),
    @('MuSCs', 'Cadm1', 'Cadm3', 'Resolving-Mac', 3, 1, 1.182122, 3.546366, 0.0566733845769154, 0.0566733845769154, 1, 0.3333333333333333, 0.01457166666666667, 0.043715, 0.001118071124106636, 0.001118071124106636, 0.01722548774333333, 0.15502938969, 0.0000633648748008394999262, 0.0000633648748008395134787),
    @('Resolving-Mac', 'Cadm1', 'Cadm3', 'ECs', 3, 1, 12.87633433333333, 38.629003, 0.6173182189435098, 0.6173182189435098, 3, 1, 8.725531, 26.176593, 0.6695022935100516, 0.6695022935100517, 112.3528543918643, 1011.175689526779, 0.41329596340822, 0.4132959634082201),
    @('Resolving-Mac', 'Cadm1', 'Cadm3', 'FAPs', 3, 1, 12.87633433333333, 38.629003, 0.6173182189435098, 0.6173182189435098, 3, 1, 2.676334666666667, 8.029004, 0.2053527971574215, 0.2053527971574215, 34.46137995589022, 310.152419603012, 0.1267680229962873, 0.1267680229962873),
    @('Resolving-Mac', 'Cadm1', 'Cadm3', 'Inflammatory-Mac', 3, 1, 12.87633433333333, 38.629003, 0.6173182189435098, 0.6173182189435098, 3, 1, 0.735084, 2.205252, 0.05640234662194688, 0.05640234662194689, 9.465187347084, 85.18668612375599, 0.03481819616089474, 0.03481819616089474),
    @('Resolving-Mac', 'Cadm1', 'Cadm3', 'MuSCs', 3, 1, 12.87633433333333, 38.629003, 0.6173182189435098, 0.6173182189435098, 3, 1, 0.8813406666666667, 2.644022, 0.06762449158647323, 0.06762449158647324, 11.34843708556289, 102.135933770066, 0.04174583070312202, 0.04174583070312202),
    @('Resolving-Mac', 'Cadm1', 'Cadm3', 'Resolving-Mac', 3, 1, 12.87633433333333, 38.629003, 0.6173182189435098, 0.6173182189435098, 1, 0.3333333333333333, 0.01457166666666667, 0.043715, 0.001118071124106636, 0.001118071124106636, 0.1876296517938889, 1.688666866145, 0.0006902056749856766, 0.0006902056749856767)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $data[$i]
    for ($j = 0; $j -lt $row.Count; $j++) {
        $ws.Cells.Item(2 + $i, 1 + $j).Value = $row[$j]
    }
}

Write-Host "done"